$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Thursday row (row 5): the CS161 class entry that was merged across
# I5:N5 is removed, and the LUNCH BREAK block (previously the single
# merged N-less run O5:U5 of individually styled cells) is turned into a
# single merge N5:U5 carrying the "LUNCH BREAK" label. ---
$ws.Range("I5:N5").UnMerge()
$ws.Range("I5:N5").ClearContents()
$ws.Range("I5:N5").ClearFormats()

$ws.Range("O5:U5").ClearContents()
$ws.Range("O5:U5").ClearFormats()

$ws.Range("N5:U5").Merge()
$ws.Range("N5").Value = "LUNCH BREAK"
$ws.Range("N2").Copy()
$ws.Range("N5").PasteSpecial(-4122)  # xlPasteFormats

# --- Monday row (row 2): add the HS161 class entry merged across
# AG2:AL2, styled like the other class-entry cells (e.g. R4). ---
$ws.Range("AG2:AL2").Merge()
$ws.Range("AG2").Value = "HS161 | English Language and Communication | Dr. Rajesh N S | C003"
$ws.Range("R4").Copy()
$ws.Range("AG2").PasteSpecial(-4122)  # xlPasteFormats

$excel.CutCopyMode = 0
